$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting the existing
# quarter columns (D:K) to (F:M).
$ws.Range("D:E").Insert()

# The inserted columns copy formatting from the column to their left
# (column C) by default; re-apply the number formatting from the
# (now-shifted) former column D -- now column F -- onto the two new
# columns so the new cells share the same style as the rest of the
# data (date format on row 7/38/80, thousands-separator number format
# everywhere else).
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("F7:F102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New values for the two newly-inserted quarter columns (D = most
# recent quarter, E = the quarter before it), keyed by row number.
$newData = @{
    7 = @(43491, 43400)
    8 = @(1396700, 1404800)
    9 = @(1097200, 1109700)
    10 = @(299500, 295100)
    11 = @($null, $null)
    12 = @("NA", "NA")
    13 = @(0, 0)
    14 = @(0, 0)
    15 = @(0, 0)
    16 = @($null, $null)
    17 = @(1351400, 1363500)
    18 = @(45300, 41300)
    19 = @($null, $null)
    20 = @(1500, 5900)
    21 = @(67300, 68200)
    22 = @(9200, 9500)
    23 = @(37600, 37700)
    24 = @(6900, 9100)
    25 = @(0, 0)
    26 = @(30700, 28600)
    27 = @(30900, 28900)
    28 = @(0, 0)
    29 = @(300, "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(-1500, -5900)
    33 = @(31200, 28900)
    34 = @(0, 0)
    35 = @(31200, 28900)
    38 = @(43491, 43400)
    39 = @($null, $null)
    40 = @($null, $null)
    41 = @(117400, 165500)
    42 = @(0, 0)
    43 = @(547200, 571200)
    44 = @(845800, 792500)
    45 = @(169200, 176000)
    46 = @(1679700, 1705300)
    47 = @(96900, 86500)
    48 = @(290300, 289600)
    49 = @(1174500, 1182500)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(79200, 80300)
    53 = @(0, 0)
    54 = @(3320500, 3344300)
    55 = @($null, $null)
    56 = @($null, $null)
    57 = @(688100, 662700)
    58 = @(51100, 80300)
    59 = @(193400, 225900)
    60 = @(932700, 968900)
    61 = @(732600, 738300)
    62 = @(182500, 181300)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(1851400, 1892300)
    67 = @($null, $null)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(1480100, 1473400)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(1469100, 1452000)
    77 = @(0, 0)
    80 = @(43491, 43400)
    81 = @(31200, 28900)
    82 = @($null, $null)
    83 = @(20500, 21100)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(-123900, -5200)
    90 = @($null, $null)
    91 = @(-11800, -11900)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(131200, 138800)
    95 = @($null, $null)
    96 = @(-24800, -24700)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-56800, -58900)
    101 = @(1400, -600)
    102 = @(-48100, 74000)
}

foreach ($r in $newData.Keys) {
    $pair = $newData[$r]
    if ($pair[0] -ne $null) {
        $ws.Cells.Item($r, 4).Value = $pair[0]
    }
    if ($pair[1] -ne $null) {
        $ws.Cells.Item($r, 5).Value = $pair[1]
    }
}
